# Updated cryptos list values (price + 1h volume change) per upstream commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.901.79"
$ws.Range("E2").Value = "'  -1.28%  "
$ws.Range("D3").Value = "'1.775.61"
$ws.Range("E3").Value = "'  -1.33%  "
$ws.Range("E4").Value = "'  +0.36%  "
$ws.Range("D5").Value = "'315.12"
$ws.Range("E5").Value = "'  -0.38%  "
$ws.Range("E6").Value = "'  +0.31%  "
$ws.Range("D7").Value = "'0.5349"
$ws.Range("E7").Value = "'  -2.62%  "
$ws.Range("D8").Value = "'0.3719"
$ws.Range("E8").Value = "'  -4.08%  "
$ws.Range("D9").Value = "'0.07385"
$ws.Range("E9").Value = "'  -2.54%  "
$ws.Range("D10").Value = "'41.38"
$ws.Range("E10").Value = "'  -3.07%  "
$ws.Range("D11").Value = "'1.089"
$ws.Range("E11").Value = "'  -2.56%  "
$ws.Range("E12").Value = "'  +0.36%  "
$ws.Range("E13").Value = "'  -3.36%  "
$ws.Range("D14").Value = "'6.058"
$ws.Range("E14").Value = "'  -2.30%  "
$ws.Range("D15").Value = "'1.779.81"
$ws.Range("E15").Value = "'  -0.81%  "
$ws.Range("D16").Value = "'7.188"
$ws.Range("D17").Value = "'87.47"
$ws.Range("E17").Value = "'  -4.89%  "
$ws.Range("D18").Value = "'0.00001049"
$ws.Range("E18").Value = "'  -1.93%  "
$ws.Range("D19").Value = "'0.06444"
$ws.Range("E19").Value = "'  -0.10%  "
$ws.Range("E20").Value = "'  +0.21%  "
$ws.Range("D21").Value = "'17.31"
$ws.Range("E21").Value = "'  +0.36%  "
$ws.Range("D22").Value = "'5.875"
$ws.Range("E22").Value = "'  -1.57%  "
$ws.Range("D23").Value = "'27.958.29"
$ws.Range("E23").Value = "'  -1.11%  "
$ws.Range("D24").Value = "'11.05"
$ws.Range("E24").Value = "'  -3.96%  "
$ws.Range("E25").Value = "'  -3.08%  "
$ws.Range("D26").Value = "'156.79"
$ws.Range("E26").Value = "'  -0.67%  "
$ws.Range("D27").Value = "'20.11"
$ws.Range("E27").Value = "'  -2.71%  "
$ws.Range("D28").Value = "'1.979.46"
$ws.Range("E28").Value = "'  -1.24%  "
$ws.Range("D29").Value = "'2.268"
$ws.Range("E29").Value = "'  -5.50%  "
$ws.Range("D30").Value = "'119.74"
$ws.Range("E30").Value = "'  -3.04%  "
$ws.Range("D31").Value = "'1.099"
$ws.Range("E31").Value = "'  -2.89%  "
$ws.Range("D32").Value = "'0.1035"
$ws.Range("E32").Value = "'  +1.53%  "
$ws.Range("D33").Value = "'3.648"
$ws.Range("E33").Value = "'  -0.60%  "
$ws.Range("D34").Value = "'5.469"
$ws.Range("E34").Value = "'  -4.54%  "
$ws.Range("D35").Value = "'0.2226"
$ws.Range("D36").Value = "'0.06345"
$ws.Range("E36").Value = "'  -0.51%  "
$ws.Range("D37").Value = "'0.02258"
$ws.Range("E37").Value = "'  -2.37%  "
$ws.Range("D38").Value = "'4.951"
$ws.Range("E38").Value = "'  -1.43%  "
$ws.Range("D39").Value = "'8.377"
$ws.Range("E39").Value = "'  -5.51%  "
$ws.Range("D40").Value = "'0.6128"
$ws.Range("E40").Value = "'  -4.18%  "
$ws.Range("D41").Value = "'1.434"
$ws.Range("E41").Value = "'  +3.71%  "
$ws.Range("D42").Value = "'10.95"
$ws.Range("E42").Value = "'  -5.58%  "
$ws.Range("D43").Value = "'1.169"
$ws.Range("E43").Value = "'  +1.02%  "
$ws.Range("E44").Value = "'  +0.25%  "
$ws.Range("D45").Value = "'13.31"
$ws.Range("E45").Value = "'  -1.55%  "
$ws.Range("D46").Value = "'3.658"
$ws.Range("E46").Value = "'  -0.60%  "
$ws.Range("D47").Value = "'0.5728"
$ws.Range("E47").Value = "'  -4.25%  "
$ws.Range("D48").Value = "'125.02"
$ws.Range("E48").Value = "'  +0.88%  "
$ws.Range("D49").Value = "'1.191"
$ws.Range("E49").Value = "'  +3.94%  "
$ws.Range("D50").Value = "'1.920"
$ws.Range("E50").Value = "'  -2.69%  "
$ws.Range("D51").Value = "'0.06811"
$ws.Range("E51").Value = "'  -1.16%  "
